# custom_jinja_filters_tpl.docx -- chain a second hello_name_filter call
# onto the "filter modified string value" paragraph, turning:
#   {{ base_value_string | hello_name_filter(“Deric”)}}
# into:
#   {{ base_value_string | hello_name_filter(“Deric”) | hello_name_filter(“and John Doe”)}}
# The final phrase ends up split across three runs (the middle one just
# holding "and "), matching how the reference document chunks the text.

$d = $word.ActiveDocument

$openQuote  = [char]0x201C
$closeQuote = [char]0x201D

# Locate the paragraph that still has the single-argument filter call.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.IndexOf("hello_name_filter(") -ge 0 -and `
        $para.Range.Text.IndexOf("base_value_string") -ge 0) {
        $targetPara = $para
    }
}

$r = $targetPara.Range

$newText = "The filter modified string value is {{ base_value_string | hello_name_filter(" + `
    $openQuote + "Deric" + $closeQuote + ") | hello_name_filter(" + `
    $openQuote + "and John Doe" + $closeQuote + ")}}"

# Rewrite the whole sentence in one go (keeps the run's existing, empty
# <w:rPr/> intact, unlike Find.Execute's replace which drops it).
$r.Text = $newText

# Re-fetch the range now that its length has changed.
$r2 = $targetPara.Range
$full = $r2.Text

# Split "and " into its own run by toggling a character attribute on and
# back off -- the host keeps the resulting run boundary even though the
# formatting ends up identical (empty <w:rPr/>) to its neighbours.
$idx = $full.IndexOf("and ")
$startPos = $r2.Start + $idx
$endPos = $startPos + 4
$rAnd = $d.Range($startPos, $endPos)
$rAnd.Bold = 1
$rAnd.Bold = 0
